$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.86481184274641
$ws.Cells.Item(2, 3).Value = 1.928141250026538
$ws.Cells.Item(2, 4).Value = -18.86481184274641
$ws.Cells.Item(2, 5).Value = -18.86481184274641
$ws.Cells.Item(2, 6).Value = -18.86481184274641
$ws.Cells.Item(2, 7).Value = -18.86481184274641
$ws.Cells.Item(2, 8).Value = -18.86481184274641
$ws.Cells.Item(2, 9).Value = -18.86481184274641
$ws.Cells.Item(2, 10).Value = -18.86481184274641
$ws.Cells.Item(2, 11).Value = -18.86481184274641

$ws.Cells.Item(3, 2).Value = -18.86481184274641
$ws.Cells.Item(3, 3).Value = -18.86481184274641
$ws.Cells.Item(3, 4).Value = -18.86481184274641
$ws.Cells.Item(3, 5).Value = -18.86481184274641
$ws.Cells.Item(3, 6).Value = -18.86481184274641
$ws.Cells.Item(3, 7).Value = -18.86481184274641
$ws.Cells.Item(3, 8).Value = -18.86481184274641
$ws.Cells.Item(3, 9).Value = 0.8821370505325354
$ws.Cells.Item(3, 10).Value = -18.86481184274641
$ws.Cells.Item(3, 11).Value = -18.86481184274641

$ws.Cells.Item(4, 2).Value = -18.86481184274641
$ws.Cells.Item(4, 3).Value = 1.976185001715787
$ws.Cells.Item(4, 4).Value = 1.653239027078162
$ws.Cells.Item(4, 5).Value = -18.86481184274641
$ws.Cells.Item(4, 6).Value = 3.379410574587161
$ws.Cells.Item(4, 7).Value = -18.86481184274641
$ws.Cells.Item(4, 8).Value = 1.525832611456367
$ws.Cells.Item(4, 9).Value = -18.86481184274641
$ws.Cells.Item(4, 10).Value = -18.86481184274641
$ws.Cells.Item(4, 11).Value = -18.86481184274641

$ws.Cells.Item(5, 2).Value = -18.86481184274641
$ws.Cells.Item(5, 3).Value = 1.656712879569866
$ws.Cells.Item(5, 4).Value = -18.86481184274641
$ws.Cells.Item(5, 5).Value = -18.86481184274641
$ws.Cells.Item(5, 6).Value = -18.86481184274641
$ws.Cells.Item(5, 7).Value = 2.68424376316933
$ws.Cells.Item(5, 8).Value = -18.86481184274641
$ws.Cells.Item(5, 9).Value = -18.86481184274641
$ws.Cells.Item(5, 10).Value = -18.86481184274641
$ws.Cells.Item(5, 11).Value = -18.86481184274641

$ws.Cells.Item(6, 2).Value = -18.86481184274641
$ws.Cells.Item(6, 3).Value = -18.86481184274641
$ws.Cells.Item(6, 4).Value = -18.86481184274641
$ws.Cells.Item(6, 5).Value = -18.86481184274641
$ws.Cells.Item(6, 6).Value = -18.86481184274641
$ws.Cells.Item(6, 7).Value = -18.86481184274641
$ws.Cells.Item(6, 8).Value = -18.86481184274641
$ws.Cells.Item(6, 9).Value = -18.86481184274641
$ws.Cells.Item(6, 10).Value = -18.86481184274641
$ws.Cells.Item(6, 11).Value = -18.86481184274641

$ws.Cells.Item(7, 2).Value = 2.443136580908365
$ws.Cells.Item(7, 3).Value = -18.86481184274641
$ws.Cells.Item(7, 4).Value = -18.86481184274641
$ws.Cells.Item(7, 5).Value = -18.86481184274641
$ws.Cells.Item(7, 6).Value = -18.86481184274641
$ws.Cells.Item(7, 7).Value = -18.86481184274641
$ws.Cells.Item(7, 8).Value = -18.86481184274641
$ws.Cells.Item(7, 9).Value = -18.86481184274641
$ws.Cells.Item(7, 10).Value = -18.86481184274641
$ws.Cells.Item(7, 11).Value = -18.86481184274641

$ws.Cells.Item(8, 2).Value = -18.86481184274641
$ws.Cells.Item(8, 3).Value = -18.86481184274641
$ws.Cells.Item(8, 4).Value = -18.86481184274641
$ws.Cells.Item(8, 5).Value = 1.786685750210047
$ws.Cells.Item(8, 6).Value = -18.86481184274641
$ws.Cells.Item(8, 7).Value = -18.86481184274641
$ws.Cells.Item(8, 8).Value = -18.86481184274641
$ws.Cells.Item(8, 9).Value = -18.86481184274641
$ws.Cells.Item(8, 10).Value = -18.86481184274641
$ws.Cells.Item(8, 11).Value = -18.86481184274641

$ws.Cells.Item(9, 2).Value = 3.864110306300137
$ws.Cells.Item(9, 3).Value = -18.86481184274641
$ws.Cells.Item(9, 4).Value = -18.86481184274641
$ws.Cells.Item(9, 5).Value = -18.86481184274641
$ws.Cells.Item(9, 6).Value = -18.86481184274641
$ws.Cells.Item(9, 7).Value = -18.86481184274641
$ws.Cells.Item(9, 8).Value = -18.86481184274641
$ws.Cells.Item(9, 9).Value = -18.86481184274641
$ws.Cells.Item(9, 10).Value = -18.86481184274641
$ws.Cells.Item(9, 11).Value = -18.86481184274641

$ws.Cells.Item(10, 2).Value = -18.86481184274641
$ws.Cells.Item(10, 3).Value = -18.86481184274641
$ws.Cells.Item(10, 4).Value = -18.86481184274641
$ws.Cells.Item(10, 5).Value = -18.86481184274641
$ws.Cells.Item(10, 6).Value = -18.86481184274641
$ws.Cells.Item(10, 7).Value = -18.86481184274641
$ws.Cells.Item(10, 8).Value = -18.86481184274641
$ws.Cells.Item(10, 9).Value = 1.880003111014668
$ws.Cells.Item(10, 10).Value = -18.86481184274641
$ws.Cells.Item(10, 11).Value = 2.112885233969762

$ws.Cells.Item(11, 2).Value = -18.86481184274641
$ws.Cells.Item(11, 3).Value = -18.86481184274641
$ws.Cells.Item(11, 4).Value = -18.86481184274641
$ws.Cells.Item(11, 5).Value = 2.92572478434944
$ws.Cells.Item(11, 6).Value = -18.86481184274641
$ws.Cells.Item(11, 7).Value = 2.915162739465451
$ws.Cells.Item(11, 8).Value = -18.86481184274641
$ws.Cells.Item(11, 9).Value = -18.86481184274641
$ws.Cells.Item(11, 10).Value = -18.86481184274641
$ws.Cells.Item(11, 11).Value = 1.957305367402156

$ws.Cells.Item(12, 2).Value = -18.86481184274641
$ws.Cells.Item(12, 3).Value = -18.86481184274641
$ws.Cells.Item(12, 4).Value = -18.86481184274641
$ws.Cells.Item(12, 5).Value = -18.86481184274641
$ws.Cells.Item(12, 6).Value = -18.86481184274641
$ws.Cells.Item(12, 7).Value = -18.86481184274641
$ws.Cells.Item(12, 8).Value = -18.86481184274641
$ws.Cells.Item(12, 9).Value = -18.86481184274641
$ws.Cells.Item(12, 10).Value = -18.86481184274641
$ws.Cells.Item(12, 11).Value = -18.86481184274641

$ws.Cells.Item(13, 2).Value = -18.86481184274641
$ws.Cells.Item(13, 3).Value = -18.86481184274641
$ws.Cells.Item(13, 4).Value = -18.86481184274641
$ws.Cells.Item(13, 5).Value = 2.526472134762181
$ws.Cells.Item(13, 6).Value = -18.86481184274641
$ws.Cells.Item(13, 7).Value = -18.86481184274641
$ws.Cells.Item(13, 8).Value = -18.86481184274641
$ws.Cells.Item(13, 9).Value = -18.86481184274641
$ws.Cells.Item(13, 10).Value = 4.321925223945822
$ws.Cells.Item(13, 11).Value = 1.776771616336089

$ws.Cells.Item(14, 2).Value = -18.86481184274641
$ws.Cells.Item(14, 3).Value = -18.86481184274641
$ws.Cells.Item(14, 4).Value = 1.521038975112659
$ws.Cells.Item(14, 5).Value = -18.86481184274641
$ws.Cells.Item(14, 6).Value = -18.86481184274641
$ws.Cells.Item(14, 7).Value = -18.86481184274641
$ws.Cells.Item(14, 8).Value = -18.86481184274641
$ws.Cells.Item(14, 9).Value = -18.86481184274641
$ws.Cells.Item(14, 10).Value = -18.86481184274641
$ws.Cells.Item(14, 11).Value = 1.842625652047411

$ws.Cells.Item(15, 2).Value = -18.86481184274641
$ws.Cells.Item(15, 3).Value = -18.86481184274641
$ws.Cells.Item(15, 4).Value = 1.723464509113474
$ws.Cells.Item(15, 5).Value = -18.86481184274641
$ws.Cells.Item(15, 6).Value = -18.86481184274641
$ws.Cells.Item(15, 7).Value = -18.86481184274641
$ws.Cells.Item(15, 8).Value = -18.86481184274641
$ws.Cells.Item(15, 9).Value = -18.86481184274641
$ws.Cells.Item(15, 10).Value = -18.86481184274641
$ws.Cells.Item(15, 11).Value = -18.86481184274641

$ws.Cells.Item(16, 2).Value = -18.86481184274641
$ws.Cells.Item(16, 3).Value = -18.86481184274641
$ws.Cells.Item(16, 4).Value = -18.86481184274641
$ws.Cells.Item(16, 5).Value = -18.86481184274641
$ws.Cells.Item(16, 6).Value = -18.86481184274641
$ws.Cells.Item(16, 7).Value = -18.86481184274641
$ws.Cells.Item(16, 8).Value = -18.86481184274641
$ws.Cells.Item(16, 9).Value = -18.86481184274641
$ws.Cells.Item(16, 10).Value = -18.86481184274641
$ws.Cells.Item(16, 11).Value = -18.86481184274641

$ws.Cells.Item(17, 2).Value = -18.86481184274641
$ws.Cells.Item(17, 3).Value = 2.171308768039124
$ws.Cells.Item(17, 4).Value = 1.863820353537773
$ws.Cells.Item(17, 5).Value = -18.86481184274641
$ws.Cells.Item(17, 6).Value = -18.86481184274641
$ws.Cells.Item(17, 7).Value = -18.86481184274641
$ws.Cells.Item(17, 8).Value = 2.064085991575442
$ws.Cells.Item(17, 9).Value = 2.025102769478679
$ws.Cells.Item(17, 10).Value = -18.86481184274641
$ws.Cells.Item(17, 11).Value = -18.86481184274641

$ws.Cells.Item(18, 2).Value = -18.86481184274641
$ws.Cells.Item(18, 3).Value = -18.86481184274641
$ws.Cells.Item(18, 4).Value = -18.86481184274641
$ws.Cells.Item(18, 5).Value = -18.86481184274641
$ws.Cells.Item(18, 6).Value = -18.86481184274641
$ws.Cells.Item(18, 7).Value = -18.86481184274641
$ws.Cells.Item(18, 8).Value = 1.991359091987735
$ws.Cells.Item(18, 9).Value = 2.125018492242985
$ws.Cells.Item(18, 10).Value = -18.86481184274641
$ws.Cells.Item(18, 11).Value = -18.86481184274641

$ws.Cells.Item(19, 2).Value = -18.86481184274641
$ws.Cells.Item(19, 3).Value = -18.86481184274641
$ws.Cells.Item(19, 4).Value = 2.086018115950153
$ws.Cells.Item(19, 5).Value = -18.86481184274641
$ws.Cells.Item(19, 6).Value = -18.86481184274641
$ws.Cells.Item(19, 7).Value = -18.86481184274641
$ws.Cells.Item(19, 8).Value = 1.544870163885079
$ws.Cells.Item(19, 9).Value = 1.874293828905967
$ws.Cells.Item(19, 10).Value = -18.86481184274641
$ws.Cells.Item(19, 11).Value = -18.86481184274641

$ws.Cells.Item(20, 2).Value = -18.86481184274641
$ws.Cells.Item(20, 3).Value = 1.084252531342641
$ws.Cells.Item(20, 4).Value = 1.483813666899082
$ws.Cells.Item(20, 5).Value = -18.86481184274641
$ws.Cells.Item(20, 6).Value = 3.262054239202753
$ws.Cells.Item(20, 7).Value = -18.86481184274641
$ws.Cells.Item(20, 8).Value = 1.692233901700149
$ws.Cells.Item(20, 9).Value = 1.249501138127743
$ws.Cells.Item(20, 10).Value = -18.86481184274641
$ws.Cells.Item(20, 11).Value = 2.256360555994753

$ws.Cells.Item(21, 2).Value = -18.86481184274641
$ws.Cells.Item(21, 3).Value = 1.311645003164959
$ws.Cells.Item(21, 4).Value = -18.86481184274641
$ws.Cells.Item(21, 5).Value = 1.67337045702908
$ws.Cells.Item(21, 6).Value = -18.86481184274641
$ws.Cells.Item(21, 7).Value = 2.592014295177993
$ws.Cells.Item(21, 8).Value = 1.492260074227563
$ws.Cells.Item(21, 9).Value = -18.86481184274641
$ws.Cells.Item(21, 10).Value = -18.86481184274641
$ws.Cells.Item(21, 11).Value = -18.86481184274641

